# edit.ps1
#
# Applies the "想去人数" (want-to-go count) refresh captured in the commit
# "Update gh-pages to output generated at 456a3b4": only column F values
# change (per-event interest counters ticking up between scrape runs) on
# the sheets 展览 (index 1), 演出 (index 2) and 全部类型 (index 4).
# 本地生活 (index 3) has no changes in this run.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1) - update F column "想去人数" counts
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 37
$ws1.Range("F4").Value = 5307
$ws1.Range("F5").Value = 177
$ws1.Range("F8").Value = 8895
$ws1.Range("F10").Value = 643
$ws1.Range("F11").Value = 12
$ws1.Range("F12").Value = 2606
$ws1.Range("F13").Value = 2606
$ws1.Range("F14").Value = 6344
$ws1.Range("F15").Value = 2340
$ws1.Range("F19").Value = 2547
$ws1.Range("F22").Value = 6583
$ws1.Range("F25").Value = 154
$ws1.Range("F27").Value = 463
$ws1.Range("F28").Value = 7205
$ws1.Range("F32").Value = 42
$ws1.Range("F35").Value = 25
$ws1.Range("F40").Value = 2553
$ws1.Range("F44").Value = 1135
$ws1.Range("F46").Value = 556
$ws1.Range("F47").Value = 3569
$ws1.Range("F49").Value = 1139
$ws1.Range("F50").Value = 27

# Sheet 2: 演出 (index 2) - update F column "想去人数" counts
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value = 210
$ws2.Range("F7").Value = 92
$ws2.Range("F10").Value = 47
$ws2.Range("F17").Value = 29

# Sheet 4: 全部类型 (index 4) - update F column "想去人数" counts
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 37
$ws4.Range("F3").Value = 5307
$ws4.Range("F4").Value = 5307
$ws4.Range("F5").Value = 177
$ws4.Range("F7").Value = 8895
$ws4.Range("F9").Value = 643
$ws4.Range("F11").Value = 2606
$ws4.Range("F14").Value = 210
$ws4.Range("F15").Value = 6344
$ws4.Range("F16").Value = 92
$ws4.Range("F20").Value = 2547
$ws4.Range("F24").Value = 6583
$ws4.Range("F28").Value = 154
$ws4.Range("F30").Value = 463
$ws4.Range("F31").Value = 7205
$ws4.Range("F34").Value = 42
$ws4.Range("F44").Value = 1135
$ws4.Range("F46").Value = 3569
$ws4.Range("F49").Value = 1139
$ws4.Range("F50").Value = 29
$ws4.Range("F51").Value = 27
